$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.253.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.43%  "

$ws.Range("D3").Value = "'1.862.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.66%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "'0.7029"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.34%  "

$ws.Range("D6").Value = "'237.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.31%  "

$ws.Range("D7").Value = "'1.0000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.27%  "

$ws.Range("D8").Value = "'0.07690"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.63%  "

$ws.Range("D9").Value = "'0.3056"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.07%  "

$ws.Range("D10").Value = "'23.40"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.24%  "

$ws.Range("D11").Value = "'0.08148"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.46%  "

$ws.Range("D12").Value = "'1.869.52"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.39%  "

$ws.Range("D13").Value = "'0.7203"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.50%  "

$ws.Range("D14").Value = "'5.164"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.47%  "

$ws.Range("D15").Value = "'89.73"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.57%  "

$ws.Range("D16").Value = "'29.266.74"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.63%  "

$ws.Range("D17").Value = "'5.761"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.68%  "

$ws.Range("D18").Value = "'13.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.45%  "

$ws.Range("D19").Value = "'238.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.51%  "

$ws.Range("D20").Value = "'0.000007715"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.32%  "

$ws.Range("D21").Value = "'0.9992"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.31%  "

$ws.Range("D22").Value = "'2.112.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.83%  "

$ws.Range("E23").Value = "  -0.16%  "

$ws.Range("D24").Value = "'7.484"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.13%  "

$ws.Range("D25").Value = "'0.1489"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.79%  "

$ws.Range("D26").Value = "'162.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.90%  "

$ws.Range("D27").Value = "'9.023"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("D28").Value = "'18.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.15%  "

$ws.Range("D29").Value = "'2.011"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.69%  "

$ws.Range("D30").Value = "'1.419"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.67%  "

$ws.Range("D31").Value = "'4.447"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.20%  "

$ws.Range("D32").Value = "'1.484"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.89%  "

$ws.Range("D33").Value = "'4.015"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.83%  "

$ws.Range("D34").Value = "'0.05210"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.59%  "

$ws.Range("E35").Value = "  -1.53%  "

$ws.Range("D36").Value = "'0.7134"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.30%  "

$ws.Range("D37").Value = "'1.001"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.02%  "

$ws.Range("D38").Value = "'2.660"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.10%  "

$ws.Range("D39").Value = "'0.01856"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.32%  "

$ws.Range("D40").Value = "'2.726"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.99%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.9365"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.68%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "'1.153.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +10.34%  "

$ws.Range("D43").Value = "'0.4299"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.44%  "

$ws.Range("D44").Value = "'71.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.70%  "

$ws.Range("D45").Value = "'5.872"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.29%  "

$ws.Range("E46").Value = "  -0.24%  "

$ws.Range("D47").Value = "'103.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.95%  "

$ws.Range("E48").Value = "  +3.09%  "

$ws.Range("D49").Value = "'2.011.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.21%  "

$ws.Range("D50").Value = "'9.159"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.18%  "

$ws.Range("D51").Value = "'6.999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.18%  "
